$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q4 at the top of
#    the data (row 2), pushing all existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# New row 2 starts out with no explicit number format; clear any
# inherited formatting on B2:D2 so they match the plain (unstyled) data
# cells used elsewhere in the column.
$summary.Range("B2:D2").ClearFormats()

# Column A carries the bold/bordered "index" style used throughout the
# sheet - copy it from the row below (A3) instead of re-building the
# font/border by hand.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.07000000000000001

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet: duplicate the "2022-Q3" tab (so headers,
#    column widths and styles match exactly) right before it, rename the
#    copy, then overwrite its data with the Q4 fund holdings.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The source tab only has 3 data rows (2-4); the Q4 tab needs a 4th
# (row 5). Give A5 the same styled look as the rows above it before
# writing the index value.
$q4.Range("A4").Copy()
$q4.Range("A5").PasteSpecial(-4122)

$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1
$q4.Range("A4").Value = 2
$q4.Range("A5").Value = 3

# Columns B:G hold numeric-looking strings ("6.76", "0.0480", fund
# codes with leading zeros, ...) that must stay text. Force a text
# number format before assigning so Excel doesn't coerce them to
# numbers, then clear the formatting back off afterwards (the stored
# cell type remains text either way) so the cells end up unstyled, like
# the rest of the workbook's data rows.
$q4.Range("B2:G5").NumberFormat = "@"

$q4.Range("B2").Value = "015784"
$q4.Range("C2").Value = "中信建投中证1000指数增强A"
$q4.Range("D2").Value = "6.76"
$q4.Range("E2").Value = "89.78"
$q4.Range("F2").Value = "0.71"
$q4.Range("G2").Value = "0.0480"
$q4.Range("H2").Value = 1

$q4.Range("B3").Value = "015785"
$q4.Range("C3").Value = "中信建投中证1000指数增强C"
$q4.Range("D3").Value = "2.40"
$q4.Range("E3").Value = "89.78"
$q4.Range("F3").Value = "0.71"
$q4.Range("G3").Value = "0.0170"
$q4.Range("H3").Value = 1

$q4.Range("B4").Value = "970046"
$q4.Range("C4").Value = "东海证券海睿健行灵活配置混合A"
$q4.Range("D4").Value = "0.14"
$q4.Range("E4").Value = "87.16"
$q4.Range("F4").Value = "4.06"
$q4.Range("G4").Value = "0.0057"
$q4.Range("H4").Value = 7

$q4.Range("B5").Value = "970047"
$q4.Range("C5").Value = "东海证券海睿健行灵活配置混合B"
$q4.Range("D5").Value = "0.09"
$q4.Range("E5").Value = "87.16"
$q4.Range("F5").Value = "4.06"
$q4.Range("G5").Value = "0.0037"
$q4.Range("H5").Value = 7

$q4.Range("B2:G5").ClearFormats()
